$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 198, shifting existing rows 198:282 down to 199:283.
$ws.Rows.Item(198).Insert()

# Populate the newly-inserted row 198 with the new data record.
$ws.Range("A198").Value = 10
$ws.Range("B198").Value = "Vega Modelo de Temuco"
$ws.Range("C198").Value = "La Araucanía"
$ws.Range("D198").Value = 44553
$ws.Range("E198").Value = 9
$ws.Range("F198").Value = 100112040
$ws.Range("G198").Value = "Cilantro"
$ws.Range("H198").Value = "Sin especificar"
$ws.Range("I198").Value = "Primera"
$ws.Range("J198").Value = 125
$ws.Range("K198").Value = 6000
$ws.Range("L198").Value = 6000
$ws.Range("M198").Value = 6000
$ws.Range("N198").Value = "$/docena de atados (2 kilos)"
$ws.Range("O198").Value = "Provincia de Cautín"
$ws.Range("P198").Value = 3000
$ws.Range("Q198").Value = 2
$ws.Range("R198").Value = "Hortaliza"
